# grade_example.xlsx — "fixed some bugs and add some features"
#
# Content edits applied to Sheet1:
#   1. Student ID column (A2:A6) renumbered from the 5121410803x series
#      to a simple 1000000000x series.
#   2. Selection moved from the old A7:J20 remnant onto the actual data
#      block A2:A6 (active cell A2).
#   3. Sheet-wide font switched from the CJK "等线" face to "Calibri".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Renumber the student IDs in column A ---------------------------
$ws.Range("A2").Value = 10000000001
$ws.Range("A3").Value = 10000000002
$ws.Range("A4").Value = 10000000003
$ws.Range("A5").Value = 10000000004
$ws.Range("A6").Value = 10000000005

# --- 2. Update the saved selection to the data range --------------------
$ws.Range("A2:A6").Select() | Out-Null

# --- 3. Switch the sheet's font to Calibri ------------------------------
$ws.Cells.Font.Name = "Calibri"
